$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: relocate the hidden "_GoBack" bookmark.
#
# It currently sits at the very end of the paragraph that reads
# "Mounting for stabilizer " (right before that paragraph's mark). The
# target state has it removed from there and instead placed at the very
# start of the following (empty) paragraph - i.e. the bookmark hops over
# the paragraph mark, the paragraph count does not change.
# ---------------------------------------------------------------------------
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -match "Mounting for stabilizer") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne $null -and $targetIndex -lt $d.Paragraphs.Count) {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $nextPara = $d.Paragraphs($targetIndex + 1)
        $destRange = $nextPara.Range
        $destRange.Collapse(1)   # wdCollapseStart -> zero-length range at the paragraph's start
        $d.Bookmarks.Add("_GoBack", $destRange)
    }
}

# ---------------------------------------------------------------------------
# Part 2: point the Discord invite link at the correct landing page.
#
# The link text/URL lives in the page header as a hyperlink run styled with
# the "Hyperlink" character style. Only the differing tail of the invite
# code is rewritten (instead of replacing the whole run's text) so the
# existing run formatting (rStyle etc.) is left untouched, matching the
# target diff which only touches the <w:t> content.
# ---------------------------------------------------------------------------
$oldUrl = "https://discord.gg/DkzJaFQWHf"
$newUrl = "https://discord.gg/ab7DeAkMmw"

$maxLen = [Math]::Min($oldUrl.Length, $newUrl.Length)
$prefixLen = 0
while ($prefixLen -lt $maxLen -and $oldUrl.Substring($prefixLen, 1) -eq $newUrl.Substring($prefixLen, 1)) {
    $prefixLen = $prefixLen + 1
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    for ($hk = 1; $hk -le 3; $hk++) {
        $hdr = $d.Sections($s).Headers($hk)
        if ($hdr.Exists) {
            $hdrRange = $hdr.Range
            for ($j = 1; $j -le $hdrRange.Hyperlinks.Count; $j++) {
                $hl = $hdrRange.Hyperlinks($j)
                if ($hl.Range.Text -eq $oldUrl) {
                    $tailRange = $hl.Range.Duplicate()
                    $tailRange.Start = $tailRange.Start + $prefixLen
                    $tailRange.Text = $newUrl.Substring($prefixLen)
                }
            }
        }
    }
}
